$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dummy")
$ws.Range("A1").Value = "DE"
$ws.Range("A2").Value = "newstring1"
$ws.Range("A3").Value = "DE"
Write-Output "done1"
